$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.604.72'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.353.58'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.40'
$ws.Range('E5').Value = '  +5.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '565.38'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').Value = '3.345.28'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.187'
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.594'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.17'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.74'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '3.872.69'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '613.47'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.25'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = '66.585.94'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = '3.321.56'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.23'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.920'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.53'
$ws.Range('E23').Value = '  +10.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.16'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.83'
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.03'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.78'
$ws.Range('E27').Value = '  +2.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.83'
$ws.Range('E28').Value = '  +5.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.75'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.66'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.83'
$ws.Range('E31').Value = '  +9.05%  '
$ws.Range('E32').Value = '  +8.87%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '565.95'
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.19'
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('D36').Value = '3.744.32'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '57.50'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').Value = '0.0₃0735'
$ws.Range('E39').Value = '  +2.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.58'
$ws.Range('E40').Value = '  +7.84%  '
$ws.Range('E41').Value = '  -2.06%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.132'
$ws.Range('E42').Value = '  +3.75%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.46'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.74'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.346'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0428'
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.25'
$ws.Range('E47').Value = '  +2.56%  '
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.63'
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('E51').Value = '  +2.39%  '
